$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "30.105.95"
$ws.Cells.Item(3, 4).Value = "1.943.37"
$ws.Cells.Item(3, 5).Value = "  +3.76%  "
$ws.Cells.Item(4, 4).Value = "0.9996"
$ws.Cells.Item(4, 5).Value = "  -0.19%  "
$ws.Cells.Item(5, 4).Value = "320.47"
$ws.Cells.Item(5, 5).Value = "  +0.34%  "
$ws.Cells.Item(6, 4).Value = "0.9991"
$ws.Cells.Item(7, 4).Value = "0.5090"
$ws.Cells.Item(7, 5).Value = "  +0.94%  "
$ws.Cells.Item(8, 4).Value = "0.4052"
$ws.Cells.Item(8, 5).Value = "  +2.41%  "
$ws.Cells.Item(9, 4).Value = "0.08415"
$ws.Cells.Item(9, 5).Value = "  +2.56%  "
$ws.Cells.Item(10, 4).Value = "1.123"
$ws.Cells.Item(10, 5).Value = "  +2.81%  "
$ws.Cells.Item(11, 4).Value = "42.29"
$ws.Cells.Item(11, 5).Value = "  +0.25%  "
$ws.Cells.Item(12, 4).Value = "24.05"
$ws.Cells.Item(12, 5).Value = "  +1.27%  "
$ws.Cells.Item(13, 4).Value = "6.426"
$ws.Cells.Item(13, 5).Value = "  +2.06%  "
$ws.Cells.Item(14, 4).Value = "1.939.68"
$ws.Cells.Item(14, 5).Value = "  +3.77%  "
$ws.Cells.Item(15, 4).Value = "7.320"
$ws.Cells.Item(15, 5).Value = "  +1.84%  "
$ws.Cells.Item(16, 4).Value = "1.001"
$ws.Cells.Item(16, 5).Value = "  -0.08%  "
$ws.Cells.Item(17, 4).Value = "93.19"
$ws.Cells.Item(17, 5).Value = "  +1.40%  "
$ws.Cells.Item(18, 4).Value = "0.00001100"
$ws.Cells.Item(18, 5).Value = "  +0.98%  "
$ws.Cells.Item(19, 4).Value = "0.06514"
$ws.Cells.Item(19, 5).Value = "  +1.68%  "
$ws.Cells.Item(20, 4).Value = "18.57"
$ws.Cells.Item(20, 5).Value = "  +2.39%  "
$ws.Cells.Item(21, 4).Value = "0.9992"
$ws.Cells.Item(21, 5).Value = "  -0.15%  "
$ws.Cells.Item(22, 4).Value = "5.991"
$ws.Cells.Item(22, 5).Value = "  +2.45%  "
$ws.Cells.Item(23, 4).Value = "30.105.24"
$ws.Cells.Item(23, 5).Value = "  +0.23%  "
$ws.Cells.Item(24, 4).Value = "11.39"
$ws.Cells.Item(24, 5).Value = "  +2.25%  "
$ws.Cells.Item(25, 5).Value = "  +0.88%  "
$ws.Cells.Item(26, 2).Value = "EthereumClassic"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(26, 4).Value = "22.34"
$ws.Cells.Item(26, 5).Value = "  +4.65%  "
$ws.Cells.Item(27, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(27, 4).Value = "2.148.36"
$ws.Cells.Item(27, 5).Value = "  +3.03%  "
$ws.Cells.Item(28, 4).Value = "162.87"
$ws.Cells.Item(28, 5).Value = "  +1.68%  "
$ws.Cells.Item(29, 4).Value = "2.328"
$ws.Cells.Item(29, 5).Value = "  +4.59%  "
$ws.Cells.Item(30, 4).Value = "129.50"
$ws.Cells.Item(30, 5).Value = "  +1.74%  "
$ws.Cells.Item(31, 4).Value = "1.138"
$ws.Cells.Item(31, 5).Value = "  +6.50%  "
$ws.Cells.Item(32, 5).Value = "  +1.26%  "
$ws.Cells.Item(33, 4).Value = "6.020"
$ws.Cells.Item(33, 5).Value = "  +1.28%  "
$ws.Cells.Item(34, 4).Value = "3.784"
$ws.Cells.Item(34, 5).Value = "  +2.77%  "
$ws.Cells.Item(35, 4).Value = "0.02461"
$ws.Cells.Item(35, 5).Value = "  +0.94%  "
$ws.Cells.Item(36, 4).Value = "5.349"
$ws.Cells.Item(36, 5).Value = "  +2.22%  "
$ws.Cells.Item(37, 4).Value = "1.263"
$ws.Cells.Item(37, 5).Value = "  +7.48%  "
$ws.Cells.Item(38, 4).Value = "0.06487"
$ws.Cells.Item(38, 5).Value = "  +1.94%  "
$ws.Cells.Item(39, 4).Value = "0.2160"
$ws.Cells.Item(39, 5).Value = "  +0.63%  "
$ws.Cells.Item(40, 4).Value = "0.6525"
$ws.Cells.Item(40, 5).Value = "  +3.40%  "
$ws.Cells.Item(41, 2).Value = "Aptos"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(41, 4).Value = "11.75"
$ws.Cells.Item(41, 5).Value = "  +3.97%  "
$ws.Cells.Item(42, 2).Value = "FraxShare"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(42, 4).Value = "8.720"
$ws.Cells.Item(42, 5).Value = "  +2.58%  "
$ws.Cells.Item(43, 4).Value = "1.224"
$ws.Cells.Item(43, 5).Value = "  +0.36%  "
$ws.Cells.Item(44, 2).Value = "EnergySwap"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(44, 4).Value = "13.41"
$ws.Cells.Item(44, 5).Value = "  +3.99%  "
$ws.Cells.Item(45, 2).Value = "Decentraland"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(45, 4).Value = "0.6096"
$ws.Cells.Item(45, 5).Value = "  +3.12%  "
$ws.Cells.Item(46, 4).Value = "2.187"
$ws.Cells.Item(46, 5).Value = "  +4.64%  "
$ws.Cells.Item(47, 5).Value = "  +0.02%  "
$ws.Cells.Item(48, 2).Value = "EOS"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Cells.Item(48, 4).Value = "1.214"
$ws.Cells.Item(48, 5).Value = "  +0.62%  "
$ws.Cells.Item(49, 2).Value = "Quant"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(49, 4).Value = "122.77"
$ws.Cells.Item(49, 5).Value = "  -0.12%  "
$ws.Cells.Item(50, 2).Value = "Aave"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(50, 4).Value = "78.48"
$ws.Cells.Item(50, 5).Value = "  +1.26%  "
$ws.Cells.Item(51, 2).Value = "WEMIXTOKEN"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(51, 4).Value = "1.132"
$ws.Cells.Item(51, 5).Value = "  +0.68%  "
